# Adding 3D Steel to the portfolio forecast
# Shift every timestamp in column A (rows 2-97) forward by 14 days,
# and refresh the actual-production values in column B for the rows
# that now carry the new day's solar generation curve.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps (rows 2 through 97) forward by 14 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 14
}

# Updated "Actual Production (MW)" readings for the new day.
$newValues = @{
    23 = 0
    24 = 2
    25 = 20
    26 = 57
    27 = 111
    28 = 175
    29 = 267
    30 = 377
    31 = 480
    32 = 604
    33 = 721
    34 = 837
    35 = 924
    36 = 996
    37 = 1052
    38 = 1139
    39 = 1216
    40 = 1247
    41 = 1310
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
